$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, thin border, center/top alignment)
# from the existing H1 header cell onto the two new header cells so the
# new columns match the look of the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data values for columns I ("I0") and J ("IF") across rows 2-17
$data = @(
    @(3, 4),
    @(6, 8),
    @(6, 8),
    @(5, 8),
    @(3, 7),
    @(1, 4),
    @(1, 5),
    @(1, 5),
    @(1, 4),
    @(1, 5),
    @(1, 4),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
